$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 7522
$ws.Range("F4").Value  = 3567
$ws.Range("F6").Value  = 3883
$ws.Range("F8").Value  = 94
$ws.Range("F9").Value  = 83
$ws.Range("F10").Value = 119
$ws.Range("F11").Value = 171
$ws.Range("F12").Value = 523
$ws.Range("F14").Value = 165
$ws.Range("F17").Value = 359
$ws.Range("F18").Value = 4235
$ws.Range("F19").Value = 4235
$ws.Range("F21").Value = 421
$ws.Range("F22").Value = 1037
$ws.Range("F23").Value = 543
$ws.Range("F24").Value = 1923
$ws.Range("F25").Value = 121
$ws.Range("F26").Value = 105
$ws.Range("F27").Value = 81
$ws.Range("F28").Value = 3098
$ws.Range("F29").Value = 2360
$ws.Range("F31").Value = 88
$ws.Range("F33").Value = 113
$ws.Range("F34").Value = 135
$ws.Range("F36").Value = 7
$ws.Range("F37").Value = 112
$ws.Range("F38").Value = 4463
$ws.Range("F39").Value = 521
$ws.Range("F43").Value = 856
$ws.Range("F44").Value = 247
$ws.Range("F45").Value = 13
$ws.Range("F46").Value = 1682
$ws.Range("F48").Value = 40
$ws.Range("F49").Value = 621

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 446
$ws.Range("F10").Value = 31
$ws.Range("F20").Value = 30
$ws.Range("F22").Value = 624
$ws.Range("F23").Value = 4

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 169

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 169
$ws.Range("F4").Value  = 7522
$ws.Range("F5").Value  = 3567
$ws.Range("F6").Value  = 3883
$ws.Range("F8").Value  = 94
$ws.Range("F9").Value  = 83
$ws.Range("F10").Value = 119
$ws.Range("F12").Value = 171
$ws.Range("F13").Value = 523
$ws.Range("F15").Value = 165
$ws.Range("F17").Value = 359
$ws.Range("F18").Value = 4235
$ws.Range("F19").Value = 4235
$ws.Range("F20").Value = 31
$ws.Range("F23").Value = 421
$ws.Range("F24").Value = 1037
$ws.Range("F25").Value = 543
$ws.Range("F26").Value = 1923
$ws.Range("F27").Value = 121
$ws.Range("F28").Value = 105
$ws.Range("F29").Value = 3098
$ws.Range("F30").Value = 2360
$ws.Range("F32").Value = 88
$ws.Range("F34").Value = 113
$ws.Range("F35").Value = 135
$ws.Range("F37").Value = 112
$ws.Range("F39").Value = 4463
$ws.Range("F41").Value = 521
$ws.Range("F45").Value = 856
$ws.Range("F46").Value = 247
$ws.Range("F47").Value = 1682
$ws.Range("F49").Value = 621
